# Fix typos in the Email column (data pulled from the database) and
# leave the active selection where the user last clicked.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "kaitkaiser@example.com"
$ws.Range("E6").Value = "bachmanuel@example.net"

$ws.Range("E7").Select()
